$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Results")
$ws.Activate()

# Rows 2-15: columns B,C,D,E all become TRUE (Create/Read/Update/Delete tests passed)
for ($r = 2; $r -le 15; $r++) {
    $ws.Range("B${r}:E${r}").Value = $true
}

# Rows 16-23: columns C,D,E become TRUE (B/"Create Test Passed" stays FALSE)
for ($r = 16; $r -le 23; $r++) {
    $ws.Range("C${r}:E${r}").Value = $true
}

# Row 24: columns C,D become TRUE (B and E stay FALSE)
$ws.Range("C24:D24").Value = $true

# Scroll the view so row 11 is at the top of the window (best effort - not all
# hosts persist window scroll position).
try {
    $excel.ActiveWindow.ScrollRow = 11
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}

# Update the active cell / selection to match the saved workbook state
$ws.Range("M18").Select()
